$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.795.79'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '2.086.61'
$ws.Range('E3').Value = '  +0.24%  '
$ws.Range('E4').Value = '  +0.02%  '
$cell = $ws.Range('D5')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '233.71'
$cell.Style = $origStyle
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('E6').Value = '  -0.25%  '
$ws.Range('E7').Value = '  -0.01%  '
$cell = $ws.Range('D8')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '58.11'
$cell.Style = $origStyle
$ws.Range('E8').Value = '  -0.97%  '
$cell = $ws.Range('D9')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.393'
$cell.Style = $origStyle
$ws.Range('E9').Value = '  +0.59%  '
$ws.Range('E10').Value = '  -0.55%  '
$ws.Range('E11').Value = '  +2.51%  '
$cell = $ws.Range('D12')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '15.27'
$cell.Style = $origStyle
$ws.Range('E12').Value = '  +2.84%  '
$ws.Range('D13').Value = '2.394.14'
$ws.Range('E13').Value = '  +0.20%  '
$cell = $ws.Range('D14')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '21.18'
$cell.Style = $origStyle
$ws.Range('E14').Value = '  +0.56%  '
$cell = $ws.Range('D15')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.778'
$cell.Style = $origStyle
$ws.Range('E15').Value = '  +0.71%  '
$cell = $ws.Range('D16')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.35'
$cell.Style = $origStyle
$ws.Range('E16').Value = '  +0.88%  '
$ws.Range('D17').Value = '2.089.97'
$ws.Range('E17').Value = '  +0.85%  '
$ws.Range('D18').Value = '37.727.98'
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('E19').Value = '  -0.51%  '
$cell = $ws.Range('D20')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '70.96'
$cell.Style = $origStyle
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('D21').Value = '0.0₃0835'
$ws.Range('E21').Value = '  +0.19%  '
$cell = $ws.Range('D22')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '229.71'
$cell.Style = $origStyle
$ws.Range('E22').Value = '  +0.50%  '
$cell = $ws.Range('D23')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = $origStyle
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  -0.31%  '
$cell = $ws.Range('D25')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.39'
$cell.Style = $origStyle
$cell = $ws.Range('D26')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '9.69'
$cell.Style = $origStyle
$ws.Range('E26').Value = '  +7.75%  '
$cell = $ws.Range('D27')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '171.48'
$cell.Style = $origStyle
$ws.Range('E27').Value = '  +1.31%  '
$ws.Range('E28').Value = '  -2.36%  '
$cell = $ws.Range('D29')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '19.48'
$cell.Style = $origStyle
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('E30').Value = '  -0.64%  '
$ws.Range('E31').Value = '  +0.31%  '
$cell = $ws.Range('D32')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.68'
$cell.Style = $origStyle
$ws.Range('E32').Value = '  -0.08%  '
$cell = $ws.Range('D33')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0631'
$cell.Style = $origStyle
$ws.Range('E33').Value = '  +0.01%  '
$cell = $ws.Range('D34')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.63'
$cell.Style = $origStyle
$ws.Range('E34').Value = '  -0.76%  '
$ws.Range('E35').Value = '  +0.27%  '
$ws.Range('E36').Value = '  -0.32%  '
$cell = $ws.Range('D37')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.33'
$cell.Style = $origStyle
$ws.Range('E37').Value = '  -1.64%  '
$ws.Range('E38').Value = '  -0.11%  '
$cell = $ws.Range('D39')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.38'
$cell.Style = $origStyle
$ws.Range('E39').Value = '  -0.39%  '
$ws.Range('E40').Value = '  +8.80%  '
$cell = $ws.Range('D41')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '101.15'
$cell.Style = $origStyle
$ws.Range('E41').Value = '  +2.76%  '
$ws.Range('E42').Value = '  -1.06%  '
$ws.Range('E43').Value = '  +1.91%  '
$ws.Range('E44').Value = '  +2.28%  '
$cell = $ws.Range('D45')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '16.79'
$cell.Style = $origStyle
$ws.Range('E45').Value = '  +1.73%  '
$ws.Range('D46').Value = '1.449.51'
$ws.Range('E46').Value = '  -0.74%  '
$cell = $ws.Range('D47')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.14'
$cell.Style = $origStyle
$ws.Range('E47').Value = '  -4.19%  '
$ws.Range('E48').Value = '  -0.54%  '
$ws.Range('E49').Value = '  -3.26%  '
$cell = $ws.Range('D50')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.96'
$cell.Style = $origStyle
$ws.Range('E50').Value = '  -2.16%  '
$ws.Range('D51').Value = '2.278.16'
$ws.Range('E51').Value = '  +0.17%  '
